$d = $word.ActiveDocument
$d.Content.Find.Execute("52+8=60", $true, $false, $false, $false, $false, $true, 1, $false, "99-88=11", 2) | Out-Null
$d.Content.Find.Execute("69-33=36", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=51", 2) | Out-Null
$d.Content.Find.Execute("4+59=63", $true, $false, $false, $false, $false, $true, 1, $false, "32+4=36", 2) | Out-Null
$d.Content.Find.Execute("48+23=71", $true, $false, $false, $false, $false, $true, 1, $false, "55-45=10", 2) | Out-Null
$d.Content.Find.Execute("33+50=83", $true, $false, $false, $false, $false, $true, 1, $false, "66-38=28", 2) | Out-Null
$d.Content.Find.Execute("19+27=46", $true, $false, $false, $false, $false, $true, 1, $false, "97-90=7", 2) | Out-Null
$d.Content.Find.Execute("82-80=2", $true, $false, $false, $false, $false, $true, 1, $false, "98-41=57", 2) | Out-Null
$d.Content.Find.Execute("49+9=58", $true, $false, $false, $false, $false, $true, 1, $false, "60-28=32", 2) | Out-Null
$d.Content.Find.Execute("27+7=34", $true, $false, $false, $false, $false, $true, 1, $false, "82-32=50", 2) | Out-Null
$d.Content.Find.Execute("27-26=1", $true, $false, $false, $false, $false, $true, 1, $false, "69-36=33", 2) | Out-Null
$d.Content.Find.Execute("5-3=2", $true, $false, $false, $false, $false, $true, 1, $false, "70-35=35", 2) | Out-Null
$d.Content.Find.Execute("81-36=45", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=73", 2) | Out-Null
$d.Content.Find.Execute("3+4=7", $true, $false, $false, $false, $false, $true, 1, $false, "89-35=54", 2) | Out-Null
$d.Content.Find.Execute("62+2=64", $true, $false, $false, $false, $false, $true, 1, $false, "56-32=24", 2) | Out-Null
$d.Content.Find.Execute("43+27=70", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=42", 2) | Out-Null
$d.Content.Find.Execute("76-24=52", $true, $false, $false, $false, $false, $true, 1, $false, "93-47=46", 2) | Out-Null
$d.Content.Find.Execute("99-36=63", $true, $false, $false, $false, $false, $true, 1, $false, "26+8=34", 2) | Out-Null
$d.Content.Find.Execute("28+14=42", $true, $false, $false, $false, $false, $true, 1, $false, "47-41=6", 2) | Out-Null
$d.Content.Find.Execute("69-61=8", $true, $false, $false, $false, $false, $true, 1, $false, "48+6=54", 2) | Out-Null
$d.Content.Find.Execute("55-36=19", $true, $false, $false, $false, $false, $true, 1, $false, "17+57=74", 2) | Out-Null
$d.Content.Find.Execute("99-4=95", $true, $false, $false, $false, $false, $true, 1, $false, "0+47=47", 2) | Out-Null
$d.Content.Find.Execute("25-24=1", $true, $false, $false, $false, $false, $true, 1, $false, "32-3=29", 2) | Out-Null
$d.Content.Find.Execute("6+18=24", $true, $false, $false, $false, $false, $true, 1, $false, "79-70=9", 2) | Out-Null
$d.Content.Find.Execute("74-23=51", $true, $false, $false, $false, $false, $true, 1, $false, "14+5=19", 2) | Out-Null
$d.Content.Find.Execute("94-13=81", $true, $false, $false, $false, $false, $true, 1, $false, "11+29=40", 2) | Out-Null
$d.Content.Find.Execute("98-70=28", $true, $false, $false, $false, $false, $true, 1, $false, "80-68=12", 2) | Out-Null
$d.Content.Find.Execute("55-32=23", $true, $false, $false, $false, $false, $true, 1, $false, "70-24=46", 2) | Out-Null
$d.Content.Find.Execute("17+2=19", $true, $false, $false, $false, $false, $true, 1, $false, "85-60=25", 2) | Out-Null
$d.Content.Find.Execute("83+9=92", $true, $false, $false, $false, $false, $true, 1, $false, "77-27=50", 2) | Out-Null
$d.Content.Find.Execute("93-71=22", $true, $false, $false, $false, $false, $true, 1, $false, "77+18=95", 2) | Out-Null
$d.Content.Find.Execute("52+28=80", $true, $false, $false, $false, $false, $true, 1, $false, "85+14=99", 2) | Out-Null
$d.Content.Find.Execute("77-12=65", $true, $false, $false, $false, $false, $true, 1, $false, "49+16=65", 2) | Out-Null
$d.Content.Find.Execute("69-34=35", $true, $false, $false, $false, $false, $true, 1, $false, "88-84=4", 2) | Out-Null
$d.Content.Find.Execute("29+4=33", $true, $false, $false, $false, $false, $true, 1, $false, "2+61=63", 2) | Out-Null
$d.Content.Find.Execute("48+48=96", $true, $false, $false, $false, $false, $true, 1, $false, "57+40=97", 2) | Out-Null
$d.Content.Find.Execute("50+32=82", $true, $false, $false, $false, $false, $true, 1, $false, "82-35=47", 2) | Out-Null
$d.Content.Find.Execute("37+45=82", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("99-66=33", $true, $false, $false, $false, $false, $true, 1, $false, "61-55=6", 2) | Out-Null
$d.Content.Find.Execute("17+60=77", $true, $false, $false, $false, $false, $true, 1, $false, "23+6=29", 2) | Out-Null
$d.Content.Find.Execute("57+5=62", $true, $false, $false, $false, $false, $true, 1, $false, "61-30=31", 2) | Out-Null
$d.Content.Find.Execute("10+50=60", $true, $false, $false, $false, $false, $true, 1, $false, "42+26=68", 2) | Out-Null
$d.Content.Find.Execute("82-81=1", $true, $false, $false, $false, $false, $true, 1, $false, "69-30=39", 2) | Out-Null
$d.Content.Find.Execute("75+9=84", $true, $false, $false, $false, $false, $true, 1, $false, "73+1=74", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $false, $false, $false, $false, $true, 1, $false, "98-10=88", 2) | Out-Null
$d.Content.Find.Execute("24+22=46", $true, $false, $false, $false, $false, $true, 1, $false, "4-2=2", 2) | Out-Null
$d.Content.Find.Execute("20-2=18", $true, $false, $false, $false, $false, $true, 1, $false, "68+0=68", 2) | Out-Null
$d.Content.Find.Execute("82+9=91", $true, $false, $false, $false, $false, $true, 1, $false, "76+4=80", 2) | Out-Null
$d.Content.Find.Execute("44+41=85", $true, $false, $false, $false, $false, $true, 1, $false, "56-19=37", 2) | Out-Null
$d.Content.Find.Execute("18+64=82", $true, $false, $false, $false, $false, $true, 1, $false, "67-44=23", 2) | Out-Null
$d.Content.Find.Execute("84-49=35", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=79", 2) | Out-Null
$d.Content.Find.Execute("42+22=64", $true, $false, $false, $false, $false, $true, 1, $false, "77-59=18", 2) | Out-Null
$d.Content.Find.Execute("48+11=59", $true, $false, $false, $false, $false, $true, 1, $false, "99-48=51", 2) | Out-Null
$d.Content.Find.Execute("87-38=49", $true, $false, $false, $false, $false, $true, 1, $false, "37-11=26", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $true, $false, $false, $false, $false, $true, 1, $false, "27+8=35", 2) | Out-Null
$d.Content.Find.Execute("41+3=44", $true, $false, $false, $false, $false, $true, 1, $false, "84-81=3", 2) | Out-Null
$d.Content.Find.Execute("26+37=63", $true, $false, $false, $false, $false, $true, 1, $false, "94-70=24", 2) | Out-Null
$d.Content.Find.Execute("54-8=46", $true, $false, $false, $false, $false, $true, 1, $false, "57+6=63", 2) | Out-Null
$d.Content.Find.Execute("62-10=52", $true, $false, $false, $false, $false, $true, 1, $false, "34-32=2", 2) | Out-Null
$d.Content.Find.Execute("85-81=4", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=18", 2) | Out-Null
$d.Content.Find.Execute("51+20=71", $true, $false, $false, $false, $false, $true, 1, $false, "87-42=45", 2) | Out-Null
$d.Content.Find.Execute("81+8=89", $true, $false, $false, $false, $false, $true, 1, $false, "20-0=20", 2) | Out-Null
$d.Content.Find.Execute("40+23=63", $true, $false, $false, $false, $false, $true, 1, $false, "68+27=95", 2) | Out-Null
$d.Content.Find.Execute("8+23=31", $true, $false, $false, $false, $false, $true, 1, $false, "18+36=54", 2) | Out-Null
$d.Content.Find.Execute("53-24=29", $true, $false, $false, $false, $false, $true, 1, $false, "2+49=51", 2) | Out-Null
$d.Content.Find.Execute("96-17=79", $true, $false, $false, $false, $false, $true, 1, $false, "19+0=19", 2) | Out-Null
$d.Content.Find.Execute("53-21=32", $true, $false, $false, $false, $false, $true, 1, $false, "75-13=62", 2) | Out-Null
$d.Content.Find.Execute("54+29=83", $true, $false, $false, $false, $false, $true, 1, $false, "57+7=64", 2) | Out-Null
$d.Content.Find.Execute("43+12=55", $true, $false, $false, $false, $false, $true, 1, $false, "48-42=6", 2) | Out-Null
$d.Content.Find.Execute("45-10=35", $true, $false, $false, $false, $false, $true, 1, $false, "30+46=76", 2) | Out-Null
$d.Content.Find.Execute("51-20=31", $true, $false, $false, $false, $false, $true, 1, $false, "40+3=43", 2) | Out-Null
$d.Content.Find.Execute("71-25=46", $true, $false, $false, $false, $false, $true, 1, $false, "27+70=97", 2) | Out-Null
$d.Content.Find.Execute("37+10=47", $true, $false, $false, $false, $false, $true, 1, $false, "50+5=55", 2) | Out-Null
$d.Content.Find.Execute("97-30=67", $true, $false, $false, $false, $false, $true, 1, $false, "95-51=44", 2) | Out-Null
$d.Content.Find.Execute("17+16=33", $true, $false, $false, $false, $false, $true, 1, $false, "13+40=53", 2) | Out-Null
$d.Content.Find.Execute("12+61=73", $true, $false, $false, $false, $false, $true, 1, $false, "45-41=4", 2) | Out-Null
$d.Content.Find.Execute("50+33=83", $true, $false, $false, $false, $false, $true, 1, $false, "56-32=24", 2) | Out-Null
$d.Content.Find.Execute("10+53=63", $true, $false, $false, $false, $false, $true, 1, $false, "90-35=55", 2) | Out-Null
$d.Content.Find.Execute("27-2=25", $true, $false, $false, $false, $false, $true, 1, $false, "59-14=45", 2) | Out-Null
$d.Content.Find.Execute("51-34=17", $true, $false, $false, $false, $false, $true, 1, $false, "85-29=56", 2) | Out-Null
$d.Content.Find.Execute("25+36=61", $true, $false, $false, $false, $false, $true, 1, $false, "59-28=31", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "54+26=80", 2) | Out-Null
$d.Content.Find.Execute("86-34=52", $true, $false, $false, $false, $false, $true, 1, $false, "82-39=43", 2) | Out-Null
$d.Content.Find.Execute("14+58=72", $true, $false, $false, $false, $false, $true, 1, $false, "46-23=23", 2) | Out-Null
$d.Content.Find.Execute("29+65=94", $true, $false, $false, $false, $false, $true, 1, $false, "4+13=17", 2) | Out-Null
$d.Content.Find.Execute("12+87=99", $true, $false, $false, $false, $false, $true, 1, $false, "84-81=3", 2) | Out-Null
$d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $true, 1, $false, "71+13=84", 2) | Out-Null
$d.Content.Find.Execute("4+30=34", $true, $false, $false, $false, $false, $true, 1, $false, "43+56=99", 2) | Out-Null
$d.Content.Find.Execute("79-34=45", $true, $false, $false, $false, $false, $true, 1, $false, "94-37=57", 2) | Out-Null
$d.Content.Find.Execute("63+25=88", $true, $false, $false, $false, $false, $true, 1, $false, "15+29=44", 2) | Out-Null
$d.Content.Find.Execute("67-31=36", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=46", 2) | Out-Null
$d.Content.Find.Execute("17+51=68", $true, $false, $false, $false, $false, $true, 1, $false, "33-25=8", 2) | Out-Null
$d.Content.Find.Execute("72-64=8", $true, $false, $false, $false, $false, $true, 1, $false, "28-19=9", 2) | Out-Null
$d.Content.Find.Execute("50+46=96", $true, $false, $false, $false, $false, $true, 1, $false, "96-45=51", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=25", 2) | Out-Null
$d.Content.Find.Execute("28-1=27", $true, $false, $false, $false, $false, $true, 1, $false, "94-75=19", 2) | Out-Null
$d.Content.Find.Execute("21+26=47", $true, $false, $false, $false, $false, $true, 1, $false, "8+55=63", 2) | Out-Null
$d.Content.Find.Execute("20-17=3", $true, $false, $false, $false, $false, $true, 1, $false, "59-59=0", 2) | Out-Null
$d.Content.Find.Execute("49+32=81", $true, $false, $false, $false, $false, $true, 1, $false, "90-30=60", 2) | Out-Null
$d.Content.Find.Execute("17+48=65", $true, $false, $false, $false, $false, $true, 1, $false, "71-52=19", 2) | Out-Null
$d.Content.Find.Execute("38+6=44", $true, $false, $false, $false, $false, $true, 1, $false, "8+2=10", 2) | Out-Null
